$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 3.272327238179451
$ws.Cells.Item(2, 3).Value = 1.626987699542094
$ws.Cells.Item(2, 4).Value = 18.71679738969934
$ws.Cells.Item(2, 5).Value = 13.86384647080068
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 37.47995879822157

# Row 3
$ws.Cells.Item(3, 2).Value = 0.01253208636536152
$ws.Cells.Item(3, 3).Value = 0.002658071450198252
$ws.Cells.Item(3, 4).Value = 0.1496068669990043
$ws.Cells.Item(3, 5).Value = 0.5333859586016987
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0.6981829834162627

# Row 4
$ws.Cells.Item(4, 2).Value = 3.272327238179451
$ws.Cells.Item(4, 3).Value = 1.626987699542094
$ws.Cells.Item(4, 4).Value = 0.1496068669990043
$ws.Cells.Item(4, 5).Value = 0.5333859586016987
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 5.582307763322248

# Row 5
$ws.Cells.Item(5, 2).Value = 3.272327238179451
$ws.Cells.Item(5, 3).Value = 1.626987699542094
$ws.Cells.Item(5, 4).Value = 0.1496068669990043
$ws.Cells.Item(5, 5).Value = 0.5333859586016987
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.582307763322248

# Row 6
$ws.Cells.Item(6, 2).Value = 3.272327238179451
$ws.Cells.Item(6, 3).Value = 1.626987699542094
$ws.Cells.Item(6, 4).Value = 0.1496068669990043
$ws.Cells.Item(6, 5).Value = 0.5333859586016987
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 5.582307763322248

# Row 7
$ws.Cells.Item(7, 2).Value = 1.445647641019636
$ws.Cells.Item(7, 3).Value = 1.626987699542094
$ws.Cells.Item(7, 4).Value = 0.1496068669990043
$ws.Cells.Item(7, 5).Value = 0.5333859586016987
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 3.755628166162433

# Row 8
$ws.Cells.Item(8, 2).Value = 3.272327238179451
$ws.Cells.Item(8, 3).Value = 1.626987699542094
$ws.Cells.Item(8, 4).Value = 0.1496068669990043
$ws.Cells.Item(8, 5).Value = 0.5333859586016987
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.582307763322248

# Row 9
$ws.Cells.Item(9, 2).Value = 1.445647641019636
$ws.Cells.Item(9, 3).Value = 208501.5462402375
$ws.Cells.Item(9, 4).Value = 0.7210945179870265
$ws.Cells.Item(9, 5).Value = 13.86384647080068
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 208517.5768288673

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1169995834814548
$ws.Cells.Item(10, 3).Value = 1.626987699542094
$ws.Cells.Item(10, 4).Value = 186123.597850132
$ws.Cells.Item(10, 5).Value = 13.86384647080068
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 186139.2056838858

# Row 11
$ws.Cells.Item(11, 2).Value = 0.1169995834814548
$ws.Cells.Item(11, 3).Value = 0.3048912486333797
$ws.Cells.Item(11, 4).Value = 0.7210945179870265
$ws.Cells.Item(11, 5).Value = 0.5333859586016987
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 1.67637130870356

# Row 12
$ws.Cells.Item(12, 2).Value = 1.445647641019636
$ws.Cells.Item(12, 3).Value = 1.626987699542094
$ws.Cells.Item(12, 4).Value = 0.1496068669990043
$ws.Cells.Item(12, 5).Value = 0.5333859586016987
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 3.755628166162433

# Row 13
$ws.Cells.Item(13, 2).Value = 3.272327238179451
$ws.Cells.Item(13, 3).Value = 1.626987699542094
$ws.Cells.Item(13, 4).Value = 0.7210945179870265
$ws.Cells.Item(13, 5).Value = 13.86384647080068
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 19.48425592650926

# Row 14
$ws.Cells.Item(14, 2).Value = 3.272327238179451
$ws.Cells.Item(14, 3).Value = 1.626987699542094
$ws.Cells.Item(14, 4).Value = 3.223369029078222
$ws.Cells.Item(14, 5).Value = 0.5333859586016987
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 8.656069925401464

# Row 15
$ws.Cells.Item(15, 2).Value = 3.272327238179451
$ws.Cells.Item(15, 3).Value = 9.98352242611593
$ws.Cells.Item(15, 4).Value = 18.71679738969934
$ws.Cells.Item(15, 5).Value = 13.86384647080068
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 45.8364935247954
